# Updating the reduction entry of the library and the table so that
# reduction outcomes are (mostly) corrected to have the prefix "Rd"
# instead of "Ab". One definition (Schwartz et al., 2006 / row 23)
# keeps the "Ab" prefix despite being in the reduction section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Coded column name" values in column I from Ab_* to Rd_*
# for the "Reduction of regular opioid use" / "Opioid use rate" rows
# (rows 22-33), except row 23 (Ab_schwartz_2006) which stays as-is.
$ws.Range("I22").Value = "Rd_soyka_2008"
$ws.Range("I24").Value = "Rd_strain_1996"
$ws.Range("I25").Value = "Rd_lingA_1976"
$ws.Range("I26").Value = "Rd_woody_2008"
$ws.Range("I27").Value = "Rd_eissenberg_1997"
$ws.Range("I28").Value = "Rd_strain_1993"
$ws.Range("I29").Value = "Rd_zaks_1972"
$ws.Range("I30").Value = "Rd_strain_1999"
$ws.Range("I31").Value = "Rd_petitjean_2001"
$ws.Range("I32").Value = "`t`nRd_shufman_1994"
$ws.Range("I33").Value = "Rd_strain_1994"

# Expand the definition text for the Ling, Charuvastra, Kaim, & Klett,
# 1976 entry (row 25) to note that it is a complex definition.
$ws.Range("F25").Value = "Index of illicit morphine use ([0, 120]). Note: this is a complex definition; for details see the original paper."
